$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<are>"
$ws.Range("C2").Value = 13

# Row 3
$ws.Range("C3").Value = 14

# Row 4
$ws.Range("C4").Value = 8

# Row 5
$ws.Range("C5").Value = 13

# Row 6
$ws.Range("C6").Value = 16

# Row 7
$ws.Range("C7").Value = 9

# Row 9
$ws.Range("C9").Value = 7

# Row 10
$ws.Range("C10").Value = 5

# Row 11
$ws.Range("C11").Value = 16

# Row 12
$ws.Range("C12").Value = 6

# Row 14
$ws.Range("B14").Value = "<are>"
$ws.Range("C14").Value = 5

# Row 15
$ws.Range("C15").Value = 7

# Row 16
$ws.Range("C16").Value = 10

# Row 17
$ws.Range("C17").Value = 10

# Row 18
$ws.Range("C18").Value = 11
